$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 6 ("Analysis"): "Created four models:" -> "Created four classification
# models:"
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(3)
$tr6 = $sh6.TextFrame.TextRange

$txt6 = $tr6.Text
$old6 = "Created four models:"
$new6 = "Created four classification models:"
$idx6 = $txt6.IndexOf($old6)
$sub6 = $tr6.Characters($idx6 + 1, $old6.Length)
$sub6.Text = $new6

# ---------------------------------------------------------------------------
# Slide 7 ("Conclusions"):
#   - append " (on-time performance vs. efficiency) " to the 2nd bullet
#   - reword the 3rd bullet
#   - remove the trailing blank paragraph
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(3)
$tr7 = $sh7.TextFrame.TextRange

$q1 = [char]8220
$q2 = [char]8221

# Remove the trailing empty paragraph first -- this round-trips the whole
# TextRange.Text property, which normalizes any curly quotes back to
# straight ones, so it must happen before the curly-quote text is inserted.
$txt7 = $tr7.Text
$trimmed7 = $txt7.Substring(0, $txt7.Length - 1)
$tr7.Text = $trimmed7

# Update the 2nd bullet.
$txt7 = $tr7.Text
$old7b = "A visual examination of these factors using Tableau show a difference in performance metrics for stops associated with a Local vs. an Express Route"
$new7b = "A visual examination of these factors using Tableau show a difference in performance metrics for stops associated with a Local vs. an Express Route (on-time performance vs. efficiency) "
$idx7b = $txt7.IndexOf($old7b)
$sub7b = $tr7.Characters($idx7b + 1, $old7b.Length)
$sub7b.Text = $new7b

# Update the 3rd bullet.
$txt7 = $tr7.Text
$old7c = "GRTC might benefit from a review of their current performance metrics to better account for the differences in objects and performance (as identified in this project) between Local and Express routes (and their stops) "
$new7c = "GRTC might benefit from a review of their current performance metrics to better account for the differences in objects and performance (as identified in this project) between Local and Express routes and adopt performance metrics for each (vs. a " + $q1 + "one size fits all" + $q2 + " performance metric approach "
$idx7c = $txt7.IndexOf($old7c)
$sub7c = $tr7.Characters($idx7c + 1, $old7c.Length)
$sub7c.Text = $new7c
